$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 51 (shifts existing rows 51-107 down to 52-108)
$ws.Rows.Item(51).Insert()

# Populate the newly inserted row 51 with the new price-report record
$ws.Range("A51").Value = 7
$ws.Range("B51").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C51").Value = "Ñuble"
$ws.Range("D51").Value = 44880
$ws.Range("E51").Value = 16
$ws.Range("F51").Value = 100112021
$ws.Range("G51").Value = "Ají"
$ws.Range("H51").Value = "Americana (o)"
$ws.Range("I51").Value = "Primera"
$ws.Range("J51").Value = 60
$ws.Range("K51").Value = 25000
$ws.Range("L51").Value = 26000
$ws.Range("M51").Value = 25500
$ws.Range("N51").Value = "$/caja 25 kilos"
$ws.Range("O51").Value = "Provincia de Limarí"
$ws.Range("P51").Value = 1020
$ws.Range("Q51").Value = 25
$ws.Range("R51").Value = "Hortaliza"
